$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "reviews_count" column (E) entirely, shifting F:K left to E:J
$ws.Range("E1").EntireColumn.Delete()
